# "Add/Edit/Delete Product (first commit)"
# - Sheet1 "Asus" -> "HP", products replaced with HP laptops (2 rows, was 3)
# - Sheet2 "Dell" -> "MSI", products replaced with MSI laptops
# - Sheet3 "Apple" removed entirely

$wb = $excel.ActiveWorkbook

# --- Sheet3 "Apple" is dropped from the workbook ---
$wsApple = $wb.Worksheets.Item("Apple")
$wsApple.Delete() | Out-Null

# --- Sheet1: Asus -> HP ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "HP"

$ws1.Range("C3").Value = "Img"
$ws1.Range("D3").Value = "Tên"
$ws1.Range("E3").Value = "Hiệu"
$ws1.Range("F3").Value = "Giá"
$ws1.Range("G3").Value = "Số lượng"

$ws1.Range("B4").Value = 1
$ws1.Range("C4").Value = "/Images/hp_pavilion_14.jpg"
$ws1.Range("D4").Value = "HP Pavilion 14"
$ws1.Range("E4").Value = "HP"
$ws1.Range("F4").Value = 16000000
$ws1.Range("G4").Value = 100

$ws1.Range("B5").Value = 2
$ws1.Range("C5").Value = "/Images/hp15_ryzen.jpg"
$ws1.Range("D5").Value = "HP 15 Ryzen"
$ws1.Range("E5").Value = "HP"
$ws1.Range("F5").Value = 22000000
$ws1.Range("G5").Value = 100

# row 3 (Zenbook Vivobook AMOLED) of the old 3-row product list is gone
$ws1.Rows.Item(6).Delete() | Out-Null

$ws1.Columns.Item(3).ColumnWidth = 23

# --- Sheet2: Dell -> MSI ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "MSI"

$ws2.Range("C3").Value = "Img"
$ws2.Range("D3").Value = "Tên"
$ws2.Range("E3").Value = "Hiệu"
$ws2.Range("F3").Value = "Giá"
$ws2.Range("G3").Value = "Số lượng"

$ws2.Range("B4").Value = 1
$ws2.Range("C4").Value = "/Images/msi_gf63.jpg"
$ws2.Range("D4").Value = "MSI GF63"
$ws2.Range("E4").Value = "MSI"
$ws2.Range("F4").Value = 20000000
$ws2.Range("G4").Value = 100

$ws2.Range("B5").Value = 2
$ws2.Range("C5").Value = "/Images/msi_pulse_gl66.jpg"
$ws2.Range("D5").Value = "MSI Pulse GL66"
$ws2.Range("E5").Value = "MSI"
$ws2.Range("F5").Value = 17000000
$ws2.Range("G5").Value = 100

$ws2.Columns.Item(3).ColumnWidth = 23

# --- selection / active sheet bookkeeping to match the saved view state ---
$ws1.Activate() | Out-Null
$ws1.Range("F10").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("G9").Select() | Out-Null
